$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("K1").Value = "ImageName"
$ws.Range("K2").Value = "QA789"

$ws.Columns.Item(11).AutoFit() | Out-Null
$ws.Columns.Item(11).ColumnWidth = 10.83

$ws.Range("K2").Select()
